# daily auto push: 2026-02-18 03:16 UTC
# Insert a new data row (2026/02/18, 水, 7, 24) right before the existing
# row 811 (2026/12/29 ...), shifting all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 811..852 down to 812..853 by inserting a whole row at 811.
$ws.Rows.Item(811).Insert()

# Populate the newly inserted row with the new reading. Column A holds a
# plain text date string (not a real date), so force text formatting
# before assigning it to keep Excel from auto-converting it to a date
# serial number, then restore the default (unstyled) cell style so the
# new row matches the formatting of the surrounding data rows.
$ws.Cells.Item(811, 1).NumberFormat = "@"
$ws.Cells.Item(811, 1).Value = "2026/02/18"
$ws.Cells.Item(811, 1).Style = "Normal"
$ws.Cells.Item(811, 2).Value = "水"
$ws.Cells.Item(811, 3).Value = 7
$ws.Cells.Item(811, 4).Value = 24
